$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values in column D remain text (matches source formatting)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '37.484.72'
$ws.Range("E2").Value = '  +5.75%  '
$ws.Range("D3").Value = '2.054.55'
$ws.Range("E3").Value = '  +4.19%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '252.88'
$ws.Range("E5").Value = '  +3.66%  '
$ws.Range("E6").Value = '  +2.88%  '
$ws.Range("D7").Value = '66.59'
$ws.Range("E7").Value = '  +17.68%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +7.35%  '
$ws.Range("D10").Value = '59.68'
$ws.Range("E10").Value = '  +2.78%  '
$ws.Range("D11").Value = '0.0768'
$ws.Range("E11").Value = '  +5.45%  '
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("E13").Value = '  -2.85%  '
$ws.Range("D14").Value = '14.93'
$ws.Range("E14").Value = '  +5.27%  '
$ws.Range("D15").Value = '2.354.70'
$ws.Range("E15").Value = '  +4.09%  '
$ws.Range("D16").Value = '22.15'
$ws.Range("E16").Value = '  +27.19%  '
$ws.Range("E17").Value = '  +7.12%  '
$ws.Range("D18").Value = '2.057.73'
$ws.Range("E18").Value = '  +4.52%  '
$ws.Range("D19").Value = '37.296.19'
$ws.Range("E19").Value = '  +5.44%  '
$ws.Range("E20").Value = '  +3.40%  '
$ws.Range("D21").Value = '0.0₃0880'
$ws.Range("E21").Value = '  +5.15%  '
$ws.Range("D22").Value = '5.45'
$ws.Range("E22").Value = '  +6.98%  '
$ws.Range("D23").Value = '240.49'
$ws.Range("E23").Value = '  +3.76%  '
$ws.Range("E24").Value = '  +6.67%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = '2.41'
$ws.Range("E26").Value = '  +5.45%  '
$ws.Range("E27").Value = '  +9.87%  '
$ws.Range("D28").Value = '161.96'
$ws.Range("E28").Value = '  -1.07%  '
$ws.Range("D29").Value = '20.04'
$ws.Range("E29").Value = '  +5.41%  '
$ws.Range("E30").Value = '  +29.37%  '
$ws.Range("D31").Value = '5.29'
$ws.Range("E31").Value = '  +9.86%  '
$ws.Range("E32").Value = '  +3.88%  '
$ws.Range("E33").Value = '  +9.70%  '
$ws.Range("D34").Value = '4.74'
$ws.Range("E34").Value = '  +10.95%  '
$ws.Range("D35").Value = '0.0626'
$ws.Range("E35").Value = '  +6.89%  '
$ws.Range("D36").Value = '2.48'
$ws.Range("E36").Value = '  +5.26%  '
$ws.Range("E37").Value = '  +4.78%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").Value = '6.02'
$ws.Range("E39").Value = '  +17.54%  '
$ws.Range("D40").Value = '2.99'
$ws.Range("E40").Value = '  +34.42%  '
$ws.Range("D41").Value = '0.104'
$ws.Range("E41").Value = '  +18.08%  '
$ws.Range("D42").Value = '1.25'
$ws.Range("E42").Value = '  +3.54%  '
$ws.Range("E43").Value = '  +4.97%  '
$ws.Range("E44").Value = '  +6.97%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '17.20'
$ws.Range("E45").Value = '  +9.51%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0218'
$ws.Range("E46").Value = '  +4.48%  '
$ws.Range("D47").Value = '96.56'
$ws.Range("E47").Value = '  +6.61%  '
$ws.Range("E48").Value = '  +7.35%  '
$ws.Range("D49").Value = '1.423.80'
$ws.Range("E49").Value = '  +4.09%  '
$ws.Range("E50").Value = '  +2.03%  '
$ws.Range("D51").Value = '46.83'
$ws.Range("E51").Value = '  +2.36%  '
